$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.492.76"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.874.32"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "'315.77"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "'0.5088"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "'0.3902"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.08366"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "'1.103"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").Value = "'41.60"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "1.873.01"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "'7.263"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "'0.00001103"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'91.28"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "'0.06729"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "'17.70"
$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'5.920"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "28.506.77"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'2.229"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "2.089.38"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").Value = "'161.88"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'20.60"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").Value = "'2.385"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "'125.65"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "'5.766"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'3.610"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "'0.02457"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").Value = "'0.06532"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'8.809"
$ws.Range("E38").Value = "  -3.65%  "
$ws.Range("D39").Value = "'5.051"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").Value = "'1.190"
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").Value = "'1.242"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "'0.6391"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").Value = "'11.08"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").Value = "'1.007"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "'0.6002"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'13.08"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").Value = "'3.684"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "'2.005"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "'121.81"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  -10.93%  "
